# Update cryptocurrency price/volume figures on the active sheet (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.683.66"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "3.148.29"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'569.72"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "'149.63"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.145.32"
$ws.Range("E8").Value = "  +2.09%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("E10").Value = "  +4.06%  "
$ws.Range("D11").Value = "'6.16"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "'0.503"
$ws.Range("E12").Value = "  +6.44%  "
$ws.Range("D13").Value = "'0.0000260"
$ws.Range("E13").Value = "  +13.33%  "
$ws.Range("D14").Value = "'38.12"
$ws.Range("E14").Value = "  +8.74%  "
$ws.Range("D15").Value = "3.664.49"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "64.774.48"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "'7.19"
$ws.Range("E17").Value = "  +6.78%  "
$ws.Range("D18").Value = "3.149.56"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "'514.50"
$ws.Range("E20").Value = "  +7.04%  "
$ws.Range("D21").Value = "'14.92"
$ws.Range("E21").Value = "  +6.29%  "
$ws.Range("D22").Value = "'0.734"
$ws.Range("E22").Value = "  +8.26%  "
$ws.Range("E23").Value = "  +8.48%  "
$ws.Range("D24").Value = "'7.84"
$ws.Range("E24").Value = "  +4.00%  "
$ws.Range("D25").Value = "'84.99"
$ws.Range("E25").Value = "  +4.35%  "
$ws.Range("E27").Value = "  +4.33%  "
$ws.Range("D28").Value = "'8.87"
$ws.Range("E28").Value = "  +11.01%  "
$ws.Range("D29").Value = "'2.18"
$ws.Range("E29").Value = "  +6.21%  "
$ws.Range("D30").Value = "'27.86"
$ws.Range("E30").Value = "  +5.92%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("E33").Value = "  +8.28%  "
$ws.Range("E34").Value = "  +9.25%  "
$ws.Range("D35").Value = "'6.58"
$ws.Range("E35").Value = "  +5.98%  "
$ws.Range("D36").Value = "'55.61"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "'484.26"
$ws.Range("E37").Value = "  +10.25%  "
$ws.Range("D38").Value = "'0.0863"
$ws.Range("E38").Value = "  +6.55%  "
$ws.Range("D39").Value = "'0.0423"
$ws.Range("E39").Value = "  +3.82%  "
$ws.Range("D40").Value = "'2.96"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "3.114.37"
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("D42").Value = "'8.64"
$ws.Range("E42").Value = "  +5.34%  "
$ws.Range("E43").Value = "  +4.86%  "
$ws.Range("E44").Value = "  +12.48%  "
$ws.Range("D45").Value = "'2.45"
$ws.Range("E45").Value = "  +15.08%  "
$ws.Range("D46").Value = "'29.52"
$ws.Range("E46").Value = "  +4.62%  "
$ws.Range("D47").Value = "0.0₃0573"
$ws.Range("E47").Value = "  +11.07%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("E50").Value = "  +10.41%  "
$ws.Range("D51").Value = "'119.84"
$ws.Range("E51").Value = "  +1.72%  "
